$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (D) / Volume(1h) (E) text values.
# "ForceText" marks D-column values that look like plain numbers
# (e.g. "1.000", "47.30") -- Excel would otherwise coerce them to
# numeric and silently drop the formatting-significant trailing zeros,
# so those cells get NumberFormat "@" before the value is written.
$updates = @(
    @{Row=2; D="24.406.22"; E="  -0.69%  "; ForceText=$False},
    @{Row=3; D="1.654.00"; E="  -2.44%  "; ForceText=$False},
    @{Row=4; D="1.004"; E="  -0.24%  "; ForceText=$True},
    @{Row=5; D="307.76"; E="  -0.39%  "; ForceText=$True},
    @{Row=6; D="0.9997"; E="  -0.27%  "; ForceText=$True},
    @{Row=7; D="0.3626"; E="  -2.75%  "; ForceText=$True},
    @{Row=8; D="47.30"; E="  -3.47%  "; ForceText=$True},
    @{Row=9; D="0.3263"; E="  -4.40%  "; ForceText=$True},
    @{Row=10; D="1.122"; E="  -4.32%  "; ForceText=$True},
    @{Row=11; D="0.06952"; E="  -6.12%  "; ForceText=$True},
    @{Row=12; D="1.000"; E="  -0.28%  "; ForceText=$True},
    @{Row=13; D="5.917"; E="  -4.51%  "; ForceText=$True},
    @{Row=14; D="19.27"; E="  -6.74%  "; ForceText=$True},
    @{Row=15; D="6.604"; E="  -3.91%  "; ForceText=$True},
    @{Row=16; D="1.654.46"; E="  -2.43%  "; ForceText=$False},
    @{Row=17; D="0.00001041"; E="  -6.34%  "; ForceText=$True},
    @{Row=18; D="0.06512"; E="  -2.65%  "; ForceText=$True},
    @{Row=19; D="0.9988"; E="  -0.47%  "; ForceText=$True},
    @{Row=20; D="76.16"; E="  -7.93%  "; ForceText=$True},
    @{Row=21; D="5.898"; E="  -6.40%  "; ForceText=$True},
    @{Row=22; D="15.70"; E="  -7.38%  "; ForceText=$True},
    @{Row=23; D="12.58"; E="  -1.42%  "; ForceText=$True},
    @{Row=24; D="24.398.60"; E="  -0.59%  "; ForceText=$False},
    @{Row=25; D="2.456"; E="  +0.51%  "; ForceText=$True},
    @{Row=26; D="2.305"; E="  -15.68%  "; ForceText=$True},
    @{Row=27; D="146.08"; E="  -2.31%  "; ForceText=$True},
    @{Row=28; D="18.36"; E="  -8.70%  "; ForceText=$True},
    @{Row=29; D="1.839.89"; E="  -2.31%  "; ForceText=$False},
    @{Row=30; D="1.195"; E="  +3.44%  "; ForceText=$True},
    @{Row=31; D="124.10"; E="  -4.73%  "; ForceText=$True},
    @{Row=32; D="4.050"; E="  -3.95%  "; ForceText=$True},
    @{Row=33; D="5.556"; E="  -15.51%  "; ForceText=$True},
    @{Row=34; D="0.08352"; E="  -4.05%  "; ForceText=$True},
    @{Row=35; D="1.683"; E="  -4.55%  "; ForceText=$True},
    @{Row=36; D="12.31"; E="  -8.25%  "; ForceText=$True},
    @{Row=37; D="5.198"; E="  -4.41%  "; ForceText=$True},
    @{Row=38; D="0.06044"; E="  -6.54%  "; ForceText=$True},
    @{Row=39; D="0.02195"; E="  -7.02%  "; ForceText=$True},
    @{Row=40; D="1.205"; E="  -4.86%  "; ForceText=$True},
    @{Row=41; D="0.2045"; E="  -5.68%  "; ForceText=$True},
    @{Row=42; D="8.148"; E="  -7.64%  "; ForceText=$True},
    @{Row=43; D="0.9996"; E="  -0.39%  "; ForceText=$True},
    @{Row=44; D="0.5843"; E="  -8.03%  "; ForceText=$True},
    @{Row=45; D="3.731"; E="  -2.01%  "; ForceText=$True},
    @{Row=46; D="12.60"; E="  -8.26%  "; ForceText=$True},
    @{Row=47; D="0.5571"; E="  -7.64%  "; ForceText=$True},
    @{Row=48; D="121.83"; E="  -4.91%  "; ForceText=$True},
    @{Row=49; D="1.935"; E="  -7.73%  "; ForceText=$True},
    @{Row=50; D="0.06904"; E=$null; ForceText=$True},
    @{Row=51; D="73.87"; E="  -5.85%  "; ForceText=$True}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.ForceText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}